$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Several "Price" cells contain strings that look like numbers
    # (e.g. "234.60"). A plain .Value assignment lets Excel's COM layer
    # auto-coerce those into real numbers (losing the original text
    # formatting / trailing zeros). Force the Text number format first so
    # the value is stored verbatim as a string, then restore the cell's
    # original (default/"Normal") style so no stray formatting is left
    # behind.
    $c = $ws.Range($range)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "43.622.55"
$ws.Range("E2").Value = "  +0.68%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.351.93"
$ws.Range("E3").Value = "  +4.56%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.13%  "

# Row 5 - BNB
Set-TextValue "D5" "234.60"
$ws.Range("E5").Value = "  +1.77%  "

# Row 6 - XRP
Set-TextValue "D6" "0.656"
$ws.Range("E6").Value = "  +3.76%  "

# Row 7 - Solana
Set-TextValue "D7" "73.29"
$ws.Range("E7").Value = "  +14.32%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.522"
$ws.Range("E9").Value = "  +19.18%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0978"
$ws.Range("E10").Value = "  +2.56%  "

# Row 11 - Avalanche
Set-TextValue "D11" "27.26"
$ws.Range("E11").Value = "  +2.60%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +2.55%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "2.705.39"
$ws.Range("E13").Value = "  +4.69%  "

# Row 14 - Chainlink
Set-TextValue "D14" "16.64"
$ws.Range("E14").Value = "  +11.86%  "

# Row 15 - Polkadot
Set-TextValue "D15" "6.61"
$ws.Range("E15").Value = "  +9.91%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.876"
$ws.Range("E16").Value = "  +6.94%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.373.66"
$ws.Range("E17").Value = "  +5.51%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "43.566.51"
$ws.Range("E18").Value = "  +0.74%  "

# Row 19 - ShibaInu
Set-TextValue "D19" "0.0000100"
$ws.Range("E19").Value = "  +3.86%  "

# Row 20 - was Uniswap, now Litecoin
$ws.Range("B20").Value = "Litecoin"
$ws.Range("C20").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D20" "75.64"
$ws.Range("E20").Value = "  +3.60%  "

# Row 21 - was Litecoin, now Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "6.40"
$ws.Range("E21").Value = "  +5.75%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "250.61"
$ws.Range("E22").Value = "  +1.76%  "

# Row 23 - WEMIXToken
Set-TextValue "D23" "3.81"
$ws.Range("E23").Value = "  -0.72%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.03%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +2.83%  "

# Row 26 - Cosmos
Set-TextValue "D26" "10.18"
$ws.Range("E26").Value = "  +4.85%  "

# Row 27 - Toncoin
$ws.Range("E27").Value = "  -2.39%  "

# Row 28 - EthereumClassic
Set-TextValue "D28" "22.45"
$ws.Range("E28").Value = "  +4.16%  "

# Row 29 - Monero
Set-TextValue "D29" "171.89"
$ws.Range("E29").Value = "  -1.10%  "

# Row 30 - ImmutableX
Set-TextValue "D30" "1.54"
$ws.Range("E30").Value = "  +7.04%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +2.51%  "

# Row 32 - Stellar
Set-TextValue "D32" "0.130"
$ws.Range("E32").Value = "  +4.33%  "

# Row 33 - Filecoin
Set-TextValue "D33" "5.06"
$ws.Range("E33").Value = "  +2.99%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.0698"
$ws.Range("E34").Value = "  +3.12%  "

# Row 35 - InternetComputer(DFINITY)
Set-TextValue "D35" "5.07"
$ws.Range("E35").Value = "  +2.92%  "

# Row 36 - RenderToken
Set-TextValue "D36" "3.73"
$ws.Range("E36").Value = "  +2.54%  "

# Row 37 - THORChain
Set-TextValue "D37" "6.57"
$ws.Range("E37").Value = "  +3.52%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "2.42"
$ws.Range("E38").Value = "  +7.20%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +5.49%  "

# Row 40 - InjectiveProtocol
Set-TextValue "D40" "19.40"
$ws.Range("E40").Value = "  +12.88%  "

# Row 41 - BinanceUSD
$ws.Range("E41").Value = "  +0.08%  "

# Row 42 - FraxShare
Set-TextValue "D42" "8.88"
$ws.Range("E42").Value = "  +1.07%  "

# Row 43 - ARBITRUM
Set-TextValue "D43" "1.16"
$ws.Range("E43").Value = "  +9.18%  "

# Row 44 - Aave
Set-TextValue "D44" "98.88"
$ws.Range("E44").Value = "  +2.43%  "

# Row 45 - was Cronos, now TrustWalletToken
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D45" "1.21"
$ws.Range("E45").Value = "  +2.87%  "

# Row 46 - was TrustWalletToken, now Cronos
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D46" "0.0963"
$ws.Range("E46").Value = "  +3.00%  "

# Row 47 - FTXToken
Set-TextValue "D47" "4.43"
$ws.Range("E47").Value = "  -1.88%  "

# Row 48 - Algorand
$ws.Range("E48").Value = "  +13.20%  "

# Row 49 - Maker
Set-TextValue "D49" "1.437.92"
$ws.Range("E49").Value = "  +0.57%  "

# Row 50 - RocketPoolETH
Set-TextValue "D50" "2.582.70"
$ws.Range("E50").Value = "  +4.87%  "

# Row 51 - was TerraClassic, now HuobiToken
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D51" "2.77"
$ws.Range("E51").Value = "  +1.14%  "
